$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.203534245491028
$ws.Range("B1").Value = 1.837275624275208
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.886786222457886
$ws.Range("E1").Value = 1.205102920532227
